# "add almost all lineal"
# Updates the generated numeric/expression values on the follower-restrictions,
# modified-point, and b/B vector sheets to the new generator run's outputs.
#
# Many of the target values are numbers written as *text* (shared strings),
# matching how the source workbook already stores them. Assigning a
# numeric-looking string straight to Range.Value makes Excel coerce it to a
# real number, so for those cells we briefly flip the cell to Text format,
# write the literal, then restore the "Normal" style (keeps the same visual
# formatting as before while keeping the stored type as text).

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------
# Restricciones_del_follower
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $ws.Range("B2") "-4.382729079133727"
Set-TextValue $ws.Range("D2") "0.0866877650392671"
Set-TextValue $ws.Range("E2") "0.36714176804725474"
Set-TextValue $ws.Range("F2") "0.35395667048189383"

Set-TextValue $ws.Range("B3") "0.38272907913372656"
Set-TextValue $ws.Range("D3") "0.9648587319705634"
Set-TextValue $ws.Range("E3") "0"
Set-TextValue $ws.Range("F3") "0.20337825316964653"

$ws.Range("A4").Value = "-16 - 2x + y_1 + 4y_2"
Set-TextValue $ws.Range("B4") "-13.81788203899294"
Set-TextValue $ws.Range("D4") "0.9761226555169311"
Set-TextValue $ws.Range("E4") "0"
Set-TextValue $ws.Range("F4") "0.23142472857106855"

Set-TextValue $ws.Range("B5") "3.3894519012124054"
Set-TextValue $ws.Range("D5") "0.8143958706897286"
Set-TextValue $ws.Range("E5") "0.3413513866142427"
Set-TextValue $ws.Range("F5") "0.05742636168621639"

Set-TextValue $ws.Range("B6") "-8.517138863787123"
Set-TextValue $ws.Range("D6") "0.23927405565041526"
Set-TextValue $ws.Range("E6") "0"
Set-TextValue $ws.Range("F6") "0.9861186567311179"

# ---------------------------------------------------------------
# Punto_modificado
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")

Set-TextValue $ws.Range("A2") "5.875840352759835"
Set-TextValue $ws.Range("B2") "4.382729079133727"
Set-TextValue $ws.Range("C2") "2.387767396848251"

# ---------------------------------------------------------------
# Vector_bf
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_bf")

Set-TextValue $ws.Range("A2") "-1.1901413818371251"
Set-TextValue $ws.Range("A3") "-3.9044906220677245"

# ---------------------------------------------------------------
# Vector_BF
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_BF")

Set-TextValue $ws.Range("A2") "-1.7308110929139415"
Set-TextValue $ws.Range("A3") "3.025790381433012"
